$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Sheet, [string]$CellRef, [string]$Val)
    $r = $Sheet.Range($CellRef)
    $r.NumberFormat = "@"
    $r.Value = $Val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "69.985.45"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").Value = "2.517.00"
$ws.Range("E3").Value = "  -0.92%  "
Set-TextValue $ws "D4" "0.998"
$ws.Range("E4").Value = "  -0.24%  "
Set-TextValue $ws "D5" "574.70"
$ws.Range("E5").Value = "  -0.54%  "
Set-TextValue $ws "D6" "166.63"
$ws.Range("E6").Value = "  -2.19%  "
$ws.Range("E7").Value = "  -0.08%  "
$ws.Range("E8").Value = "  +1.57%  "
$ws.Range("D9").Value = "2.517.03"
$ws.Range("E9").Value = "  -0.92%  "
Set-TextValue $ws "D10" "0.161"
$ws.Range("E10").Value = "  -2.85%  "
$ws.Range("E11").Value = "  -0.75%  "
$ws.Range("E12").Value = "  +3.51%  "
Set-TextValue $ws "D13" "4.94"
$ws.Range("E13").Value = "  +2.60%  "
$ws.Range("D14").Value = "2.983.69"
$ws.Range("E14").Value = "  -0.75%  "
$ws.Range("D15").Value = "69.620.93"
$ws.Range("E15").Value = "  -1.36%  "
$ws.Range("E16").Value = "  -2.69%  "
Set-TextValue $ws "D17" "24.96"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "2.521.53"
$ws.Range("E18").Value = "  -0.47%  "
Set-TextValue $ws "D19" "11.39"
$ws.Range("E19").Value = "  -1.94%  "
$ws.Range("E20").Value = "  +4.66%  "
Set-TextValue $ws "D21" "350.88"
$ws.Range("E21").Value = "  -2.60%  "
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("E23").Value = "  +0.76%  "
$ws.Range("E24").Value = "  +0.10%  "
Set-TextValue $ws "D25" "70.38"
$ws.Range("E25").Value = "  +1.07%  "
Set-TextValue $ws "D26" "4.01"
$ws.Range("E26").Value = "  -1.76%  "
Set-TextValue $ws "D27" "8.97"
$ws.Range("E27").Value = "  -3.04%  "
$ws.Range("D28").Value = "2.646.94"
$ws.Range("E28").Value = "  -0.95%  "
Set-TextValue $ws "D29" "0.999"
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("D30").Value = "0.0₃0900"
Set-TextValue $ws "D31" "7.91"
$ws.Range("E31").Value = "  +0.38%  "
Set-TextValue $ws "D32" "466.53"
$ws.Range("E32").Value = "  -3.76%  "
Set-TextValue $ws "D33" "1.25"
$ws.Range("E33").Value = "  -1.62%  "
$ws.Range("E34").Value = "  -1.27%  "
Set-TextValue $ws "D35" "1.00"
$ws.Range("E35").Value = "  +0.06%  "
Set-TextValue $ws "D36" "157.83"
$ws.Range("E36").Value = "  +0.41%  "
$ws.Range("E37").Value = "  +0.10%  "
Set-TextValue $ws "D38" "19.03"
$ws.Range("E38").Value = "  +1.08%  "
Set-TextValue $ws "D39" "18.55"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("E41").Value = "  +0.42%  "
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  -3.41%  "
Set-TextValue $ws "D44" "38.43"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("E45").Value = "  -6.15%  "
$ws.Range("E46").Value = "  -13.42%  "
Set-TextValue $ws "D47" "142.49"
$ws.Range("E47").Value = "  -1.51%  "
$ws.Range("E48").Value = "  +0.21%  "
$ws.Range("E49").Value = "  -0.96%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D50" "0.0731"
$ws.Range("E50").Value = "  -0.68%  "
$ws.Range("B51").Value = "Optimism"
$ws.Range("C51").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
Set-TextValue $ws "D51" "1.58"
$ws.Range("E51").Value = "  -3.24%  "
